# Import the data dictionaries: tidy up the "Variables" sheet of the
# 1_0_monthly_repeated_measures dictionary.
#
#   * drop the redundant "table" column (every row repeats the same
#     dictionary name, so it carries no information once the file lives
#     in a per-table tab)
#   * drop the unused "alias" column
#   * rename the "label:en" header to "label"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Remove column A ("table"); everything shifts one column left.
$ws.Columns.Item(1).Delete()

# Remove the (now shifted) "alias" column, originally column F, now column E.
$ws.Columns.Item(5).Delete()

# The old "label:en" header (now in column D) becomes simply "label".
$ws.Range("D1").Value = "label"

# Restore a sane selection on the sheet.
[void]$ws.Range("G16").Select()
